$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
